$wb = $excel.ActiveWorkbook

# Sheet1 ("names") holds a pool of unused random ids. Consume the first one.
$namesSheet = $wb.Worksheets.Item(1)
$consumedId = $namesSheet.Cells.Item(1, 1).Value()

# Remove the consumed id from the top of the pool; remaining ids shift up.
$namesSheet.Rows.Item(1).Delete()

# Sheet2 ("used") is an append-only log of consumed ids with metadata.
$usedSheet = $wb.Worksheets.Item(2)
$nextRow = $usedSheet.UsedRange.Rows.Count + 1

$usedSheet.Cells.Item($nextRow, 1).Value = $consumedId
$usedSheet.Cells.Item($nextRow, 2).Value = "ChatGPT Image 2026年1月18日 09_58_09.png"
$usedSheet.Cells.Item($nextRow, 3).Value = "2026-01-18 09:59:49"
